$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.572.25"
$ws.Range("E2").Value = "  -2.61%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.924.05"
$ws.Range("E3").Value = "  -2.75%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.88"
$ws.Range("E5").Value = "  +4.84%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.66"
$ws.Range("E6").Value = "  +0.62%  "

# Row 7
$ws.Range("E7").Value = "  -5.25%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.735"
$ws.Range("E9").Value = "  -3.92%  "

# Row 10
$ws.Range("E10").Value = "  -4.02%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.84"
$ws.Range("E11").Value = "  +13.83%  "

# Row 12
$ws.Range("E12").Value = "  -1.75%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.49"
$ws.Range("E13").Value = "  -2.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.553.58"
$ws.Range("E14").Value = "  -2.38%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.915.11"
$ws.Range("E15").Value = "  -2.69%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.91"
$ws.Range("E16").Value = "  -1.17%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.26"
$ws.Range("E17").Value = "  -3.91%  "

# Row 18
$ws.Range("E18").Value = "  -0.64%  "

# Row 19
$ws.Range("E19").Value = "  -3.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.525.05"
$ws.Range("E20").Value = "  -2.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "425.95"
$ws.Range("E21").Value = "  -4.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "95.99"
$ws.Range("E22").Value = "  -7.89%  "

# Row 23
$ws.Range("E23").Value = "  -1.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.21"
$ws.Range("E24").Value = "  +6.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.18"
$ws.Range("E25").Value = "  -2.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.20"
$ws.Range("E26").Value = "  -2.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.48"
$ws.Range("E27").Value = "  -4.67%  "

# Row 28
$ws.Range("E28").Value = "  +0.69%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.62"
$ws.Range("E29").Value = "  +16.43%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.18"
$ws.Range("E30").Value = "  -6.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.43"
$ws.Range("E31").Value = "  +9.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.34"
$ws.Range("E32").Value = "  -2.35%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.129"
$ws.Range("E33").Value = "  +1.17%  "

# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "48.36"
$ws.Range("E34").Value = "  +17.09%  "

# Row 35
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "676.51"
$ws.Range("E35").Value = "  +0.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.23"
$ws.Range("E36").Value = "  -2.59%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.429"
$ws.Range("E37").Value = "  -0.50%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0815"
$ws.Range("E38").Value = "  -5.36%  "

# Row 39
$ws.Range("E39").Value = "  -1.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  -3.32%  "

# Row 41
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.07%  "

# Row 42
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.32"
$ws.Range("E42").Value = "  +4.50%  "

# Row 43
$ws.Range("E43").Value = "  +0.26%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0479"
$ws.Range("E44").Value = "  -1.82%  "

# Row 45
$ws.Range("E45").Value = "  -5.45%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.69"
$ws.Range("E46").Value = "  -1.47%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.58"
$ws.Range("E47").Value = "  +5.72%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.35"
$ws.Range("E48").Value = "  -4.17%  "

# Row 49
$ws.Range("E49").Value = "  -3.56%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000271"
$ws.Range("E50").Value = "  +1.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.53"
$ws.Range("E51").Value = "  +1.17%  "
